$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) from N1 out to Z1, copying the bold/border/centered
# format from the existing N1 header cell, then set the shared "Unnamed: 1" text.
$ws.Range("N1").Copy()
$ws.Range("O1:Z1").PasteSpecial(-4122)
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(1, $col).Value = "Unnamed: 1" }

# Row 4 (Kalenderwoche), row 5 (Theke) and row 6 (Zweigstelle) carry the new,
# full 24-column (C:Z) breakdown of the underlying data upload.
$row4Vals = @("24. KW", "24. KW", "24. KW", "24. KW", "23. KW", "23. KW", "23. KW", "23. KW", "23. KW", "24. KW", "23. KW", "23. KW", "23. KW", "24. KW", "24. KW", "24. KW", "24. KW", "24. KW", "24. KW", "23. KW", "23. KW", "24. KW", "23. KW", "23. KW")
$col = 3
foreach ($v in $row4Vals) {
    $ws.Cells.Item(4, $col).Value = $v
    $col = $col + 1
}

$row5Vals = @("Auskunftstheke", "Ausleihtheke", "Auskunftstheke", "Ausleihtheke", "Ausleihtheke", "Auskunftstheke", "Ausleihtheke", "Auskunftstheke", "Auskunftstheke", "Ausleihtheke", "Ausleihtheke", "Auskunftstheke", "Ausleihtheke", "Auskunftstheke", "Auskunftstheke", "Ausleihtheke", "Ausleihtheke", "Ausleihtheke", "Auskunftstheke", "Auskunftstheke", "Auskunftstheke", "Auskunftstheke", "Ausleihtheke", "Ausleihtheke")
$col = 3
foreach ($v in $row5Vals) {
    $ws.Cells.Item(5, $col).Value = $v
    $col = $col + 1
}

$row6Vals = @("Deutz", "Deutz", "Deutz", "Deutz", "Deutz", "Deutz", "Deutz", "Deutz", "GM", "GM", "GM", "GM", "GM", "GM", "GM", "GM", "Südstadt", "Südstadt", "Südstadt", "Südstadt", "Südstadt", "Südstadt", "Südstadt", "Südstadt")
$col = 3
foreach ($v in $row6Vals) {
    $ws.Cells.Item(6, $col).Value = $v
    $col = $col + 1
}

for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(9, $col).Value = 8 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(10, $col).Value = 4 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(11, $col).Value = 4 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(12, $col).Value = 12 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(13, $col).Value = 5 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(15, $col).Value = 7 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(16, $col).Value = 19 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(17, $col).Value = 13 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(18, $col).Value = 1 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(19, $col).Value = 5 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(20, $col).Value = 19 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(21, $col).Value = 4 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(22, $col).Value = 15 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(23, $col).Value = 45 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(24, $col).Value = 18 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(25, $col).Value = 8 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(26, $col).Value = 11 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(27, $col).Value = 8 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(29, $col).Value = 10 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(31, $col).Value = 8 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(32, $col).Value = 2 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(34, $col).Value = 21 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(35, $col).Value = 6 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(36, $col).Value = 5 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(37, $col).Value = 10 }
for ($col = 15; $col -le 26; $col++) { $ws.Cells.Item(38, $col).Value = 12 }

# Ensure the sheet's used range reflects the new Z38 extent.
$ws.Range("A1").Select()
